$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: numeric 0, bold font + thin border all round + center/top alignment
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment = -4160     # xlTop
$ws.Range("B1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("B1").Borders.Weight = 2            # xlThin

# A2: numeric 0, same style as B1 -- copy the format so the engine reuses
# the same cell-style record instead of minting a near-duplicate one.
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)           # xlPasteFormats

# B2: shared string label, default (unstyled) cell
$ws.Range("B2").Value = "disconnected_elements"
